# ---------------------------------------------------------------------------
# Target edit (per the supplied OOXML diff):
#
#   1. In the "War test" paragraph, split the trailing run " test" into
#      " t" + "est" and relocate the existing "_GoBack" bookmark (which
#      previously sat alone in the following, otherwise-empty paragraph) to
#      the new split point between "t" and "est". The paragraph that used
#      to hold the bookmark becomes a plain empty paragraph.
#
#   2. Remove the now-unused Comment/Balloon Word styles (CommentReference,
#      CommentText, CommentTextChar, CommentSubject, CommentSubjectChar,
#      BalloonText, BalloonTextChar) from the style sheet.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Relocate the "_GoBack" bookmark into the middle of "War test" -----
# Find the anchor text "War t" (case-/whole-word-insensitive exact match);
# its end position is exactly where the bookmark needs to move to.
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("War t", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)

if ($found) {
    $splitPos = $findRange.End

    # Remove the bookmark from its old (empty) paragraph ...
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()

    # ... and re-add it, collapsed to a point, right after "War t" so it
    # now sits between the "t" and "est" runs.
    $splitRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $splitRange)
}

# --- 2. Drop the obsolete Comment/Balloon styles ---------------------------
# Delete tail-first: the live Styles collection re-indexes after every
# removal, so walking the names in reverse collection order keeps each
# lookup-by-name valid at the moment it runs.
$obsoleteStyles = @(
    "BalloonTextChar",
    "BalloonText",
    "CommentSubjectChar",
    "CommentSubject",
    "CommentTextChar",
    "CommentText",
    "CommentReference"
)

foreach ($styleName in $obsoleteStyles) {
    $style = $d.Styles($styleName)
    if ($style -ne $null) {
        $style.Delete()
    }
}
